# Excel COM-interop script implementing the commit's spreadsheet edit:
# Insert a new "Minority interest" row into the standardized Balance Sheet
# mapping table (standardized_BS.xlsx) right before the existing
# "Total Stockholders Equity" row, and refresh the sheet/window view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 38 (pushes the former rows 38-40 down to 39-41) ---
$ws.Rows.Item(38).Insert()

# --- Populate the newly inserted row with the "Minority interest" mapping ---
$ws.Range("A38").Value = "Minority interest"
$ws.Range("B38").Value = "MinorityInterest"
$ws.Range("C38").Value = "Total of all stockholders' equity (deficit) items, net of receivables from officers, directors, owners, and affiliates of the entity which is directly or indirectly attributable to that ownership interest in subsidiary equity which is not attributable to the parent (that is, noncontrolling interest, previously referred to as minority interest)."
$ws.Range("D38").Value = "NO"
$ws.Range("E38").Value = "BS"

# Match the row height/style of the surrounding rows (15pt custom height,
# default/normal cell style on D38 so it relies on the column's style).
$ws.Rows.Item(38).RowHeight = 15
$ws.Range("D38").Style = "Normal"

# --- Refresh the view: zoom in and scroll so row 37 (A7 region) is visible,
#     with the new cell A37 (old A2 in the pre-insert numbering shifted) selected ---
$win = $excel.ActiveWindow
$win.Zoom = 140
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("A37").Select()

Write-Host "Row 38 inserted: A38=" $ws.Range("A38").Value() "B38=" $ws.Range("B38").Value()
Write-Host "Dimension now covers through row" $ws.UsedRange.Rows.Count
